# Adding the changes we made on may 9th
#
# The data table (x, y, z sensor readings) gained 7 new rows inserted right
# after the header and 3 new rows appended at the end, while all of the
# previously existing rows shifted down by 7 (old row 2 is now row 9, ...,
# old row 21 is now row 28). We simply rewrite the full data block with its
# final contents so every cell ends up with the correct value without
# disturbing any cell formatting/styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.41813725233078, -0.3617849349975586, 0.7021896243095398),
    @(0.2964223623275757, 0.0287106670439243, 0.1348485052585601),
    @(-0.015118914656341, 0.0493273697793483, 0.1067487001419067),
    @(-0.0705549344420433, -0.0940732508897781, -0.0430659987032413),
    @(0.0039706239476799, 0.1137736514210701, 0.0995710343122482),
    @(0.0445931628346443, 0.0633772686123848, 0.0119118718430399),
    @(-0.1004873365163803, 0.0113010071218013, 0.012980886735022),
    @(-0.1327104717493057, 0.0224492978304624, 0.07544185966253281),
    @(-0.0704022198915481, 0.0503963828086853, -0.0980438739061355),
    @(0.0665843114256858, 0.0959058403968811, -0.0583376325666904),
    @(0.0452040284872055, 0.1351539343595504, -0.1539380401372909),
    @(0.2449569702148437, 0.4401284158229828, -0.2344195395708084),
    @(0.5047274231910706, 1.391398310661316, 1.06214189529419),
    @(0.0247400421649217, -0.4867068827152252, -1.290605545043945),
    @(0.8316930532455444, -1.329548239707947, 0.3587306141853332),
    @(0.7304421067237854, -2.097558498382568, -0.8232936263084412),
    @(0.2547308206558227, -0.4392121136188507, -0.3918700516223907),
    @(0.09239336848258969, -0.5294674634933472, 0.3912591934204101),
    @(-1.268614411354065, -1.000597238540649, -1.631621122360229),
    @(1.631926536560059, 8.462469100952148, -0.6192646622657776),
    @(-1.297935962677002, -0.3597996234893799, 0.7938193678855896),
    @(-0.845132052898407, 0.3094032406806946, -0.9367618560791016),
    @(-1.232115149497986, 0.111024759709835, -0.7185302376747131),
    @(-0.2264782935380935, -0.9043859839439392, 1.454012036323547),
    @(0.5639813542366028, 0.6252205967903137, -0.7906123399734497),
    @(-0.1081231459975242, 0.0042760567739605, 1.263422012329102),
    @(0.1950187236070633, -0.3927863538265228, 0.2987131178379059),
    @(0.1111774742603302, -0.1542434692382812, 0.3286454975605011),
    @(0.7035640478134155, -0.6233879923820496, 0.6962336897850037),
    @(0.2559525370597839, 0.1485929638147354, -0.1579086631536483)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
